# Insert two new data rows (for Fecha = 44917) at the top of the Apio / Vega
# Monumental Concepción block, pushing the existing rows 262-375 down to
# 264-377. This mirrors the structure of the other "Primera"/"Segunda"
# quality-pair rows already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 262.
$ws.Rows.Item(262).Insert()
$ws.Rows.Item(263).Insert()

# New row 262: Primera
$ws.Cells.Item(262, 1).Value = 11
$ws.Cells.Item(262, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(262, 3).Value = "Bíobío"
$ws.Cells.Item(262, 4).Value = 44917
$ws.Cells.Item(262, 5).Value = 8
$ws.Cells.Item(262, 6).Value = 100112017
$ws.Cells.Item(262, 7).Value = "Apio"
$ws.Cells.Item(262, 8).Value = "Americana (o)"
$ws.Cells.Item(262, 9).Value = "Primera"
$ws.Cells.Item(262, 10).Value = 100
$ws.Cells.Item(262, 11).Value = 8000
$ws.Cells.Item(262, 12).Value = 8500
$ws.Cells.Item(262, 13).Value = 8250
$ws.Cells.Item(262, 14).Value = "`$/docena de matas"
$ws.Cells.Item(262, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(262, 16).Value = 1375
$ws.Cells.Item(262, 17).Value = 6
$ws.Cells.Item(262, 18).Value = "Hortaliza"

# New row 263: Segunda
$ws.Cells.Item(263, 1).Value = 11
$ws.Cells.Item(263, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(263, 3).Value = "Bíobío"
$ws.Cells.Item(263, 4).Value = 44917
$ws.Cells.Item(263, 5).Value = 8
$ws.Cells.Item(263, 6).Value = 100112017
$ws.Cells.Item(263, 7).Value = "Apio"
$ws.Cells.Item(263, 8).Value = "Americana (o)"
$ws.Cells.Item(263, 9).Value = "Segunda"
$ws.Cells.Item(263, 10).Value = 50
$ws.Cells.Item(263, 11).Value = 6500
$ws.Cells.Item(263, 12).Value = 6500
$ws.Cells.Item(263, 13).Value = 6500
$ws.Cells.Item(263, 14).Value = "`$/docena de matas"
$ws.Cells.Item(263, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(263, 16).Value = 1083
$ws.Cells.Item(263, 17).Value = 6
$ws.Cells.Item(263, 18).Value = "Hortaliza"

# Apply the same number format (date) style used by the rest of column D to
# the two new D cells, matching the existing cells immediately below.
$ws.Range("D262").NumberFormat = $ws.Range("D264").NumberFormat
$ws.Range("D263").NumberFormat = $ws.Range("D264").NumberFormat
